$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 327.81818
$ws.Range("I41").Value = 412.33334
$ws.Range("J41").Value = 226.4
$ws.Range("K41").Value = 412.33334
$ws.Range("L41").Value = 226.4
$ws.Range("M41").Value = 27.66665999999998
$ws.Range("N41").Value = -1106.4

$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16498

$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -52488

$ws.Range("H93").Value = 31186.428
$ws.Range("J93").Value = 31186.428
$ws.Range("L93").Value = 31186.428
$ws.Range("N93").Value = -36178.428

$ws.Range("H127").Value = 1102.35
$ws.Range("I127").Value = 679.4
$ws.Range("J127").Value = 1243.3334
$ws.Range("K127").Value = 2038.2
$ws.Range("L127").Value = 3730.0002
$ws.Range("M127").Value = 2921.8
$ws.Range("N127").Value = -13650.0002

$ws.Range("H129").Value = 3442.3606
$ws.Range("I129").Value = 430.5
$ws.Range("J129").Value = 3896.9812
$ws.Range("K129").Value = 1291.5
$ws.Range("L129").Value = 11690.9436
$ws.Range("M129").Value = 3708.5
$ws.Range("N129").Value = -21690.9436

$ws.Range("H132").Value = 5072.525
$ws.Range("I132").Value = 4786.8687
$ws.Range("K132").Value = 14360.6061
$ws.Range("M132").Value = -11830.6061

$ws.Range("H137").Value = 8334231
$ws.Range("I137").Value = 852.9286
$ws.Range("K137").Value = 2558.7858
$ws.Range("M137").Value = -8.785799999999654

$ws.Range("H138").Value = 3500
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10500
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -20780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7756.452
$ws.Range("I32").Value = 7177.366
$ws.Range("J32").Value = 10919.154
$ws.Range("K32").Value = 7177.366
$ws.Range("L32").Value = 10919.154
$ws.Range("M32").Value = -6890.366
$ws.Range("N32").Value = -11493.154

$ws.Range("H112").Value = 34000
$ws.Range("J112").Value = 34000
$ws.Range("L112").Value = 34000
$ws.Range("N112").Value = -36954

$ws.Range("H132").Value = 7355126.5
$ws.Range("I132").Value = 8622552
$ws.Range("K132").Value = 25867656
$ws.Range("M132").Value = -25865126

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882

$ws.Range("H134").Value = 2550.628
$ws.Range("I134").Value = 1711.1666
$ws.Range("K134").Value = 5133.4998
$ws.Range("M134").Value = -2598.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16675018
$ws.Range("I31").Value = 9468
$ws.Range("J31").Value = 66671668
$ws.Range("K31").Value = 9468
$ws.Range("L31").Value = 66671668
$ws.Range("M31").Value = -9173
$ws.Range("N31").Value = -66672258

$ws.Range("H34").Value = 16675018
$ws.Range("I34").Value = 9468
$ws.Range("J34").Value = 66671668
$ws.Range("K34").Value = 9468
$ws.Range("L34").Value = 66671668
$ws.Range("M34").Value = -9266
$ws.Range("N34").Value = -66672072

$ws.Range("H132").Value = 20836498
$ws.Range("J132").Value = 3202.3635
$ws.Range("L132").Value = 9607.0905
$ws.Range("N132").Value = -14667.0905

$ws.Range("H140").Value = 49235.8
$ws.Range("J140").Value = 49235.8
$ws.Range("L140").Value = 49235.8
$ws.Range("N140").Value = -59595.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 15723
$ws.Range("I97").Value = 34519
$ws.Range("J97").Value = 1626
$ws.Range("K97").Value = 103557
$ws.Range("L97").Value = 4878
$ws.Range("M97").Value = -103061
$ws.Range("N97").Value = -5870

$ws.Range("H122").Value = 603.8214
$ws.Range("I122").Value = 589.14813
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5302.33317
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2852.33317
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34142.855
$ws.Range("I70").Value = 102500
$ws.Range("J70").Value = 6800
$ws.Range("K70").Value = 102500
$ws.Range("L70").Value = 6800
$ws.Range("M70").Value = -102230
$ws.Range("N70").Value = -7340

$ws.Range("H73").Value = 34142.855
$ws.Range("I73").Value = 102500
$ws.Range("J73").Value = 6800
$ws.Range("K73").Value = 102500
$ws.Range("L73").Value = 6800
$ws.Range("M73").Value = -101564
$ws.Range("N73").Value = -8672

$ws.Range("H111").Value = 35293
$ws.Range("J111").Value = 35293
$ws.Range("L111").Value = 35293
$ws.Range("N111").Value = -41427

$ws.Range("H113").Value = 44589
$ws.Range("I113").Value = 51136.6
$ws.Range("J113").Value = 938.3333
$ws.Range("K113").Value = 51136.6
$ws.Range("L113").Value = 938.3333
$ws.Range("M113").Value = -48966.6
$ws.Range("N113").Value = -5278.3333

$ws.Range("H126").Value = 5336.773
$ws.Range("I126").Value = 4202.4
$ws.Range("J126").Value = 5670.4116
$ws.Range("K126").Value = 12607.2
$ws.Range("L126").Value = 17011.2348
$ws.Range("M126").Value = -10137.2
$ws.Range("N126").Value = -21951.2348

$ws.Range("H132").Value = 3477.8936
$ws.Range("I132").Value = 2783.3
$ws.Range("J132").Value = 4703.647
$ws.Range("K132").Value = 8349.900000000001
$ws.Range("L132").Value = 14110.941
$ws.Range("M132").Value = -5819.900000000001
$ws.Range("N132").Value = -19170.941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5804.5454
$ws.Range("I122").Value = 6863.077
$ws.Range("K122").Value = 20589.231
$ws.Range("M122").Value = -18139.231

$ws.Range("H132").Value = 9811143
$ws.Range("I132").Value = 4692.857
$ws.Range("J132").Value = 21749432
$ws.Range("K132").Value = 14078.571
$ws.Range("L132").Value = 65248296
$ws.Range("M132").Value = -11548.571
$ws.Range("N132").Value = -65253356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 42250
$ws.Range("I69").Value = 35000
$ws.Range("J69").Value = 44666.668
$ws.Range("K69").Value = 35000
$ws.Range("L69").Value = 44666.668
$ws.Range("M69").Value = -34251
$ws.Range("N69").Value = -46164.668

$ws.Range("H72").Value = 42250
$ws.Range("I72").Value = 35000
$ws.Range("J72").Value = 44666.668
$ws.Range("K72").Value = 105000
$ws.Range("L72").Value = 134000.004
$ws.Range("M72").Value = -101256
$ws.Range("N72").Value = -141488.004

$ws.Range("H132").Value = 2258.3784
$ws.Range("I132").Value = 1821.9667
$ws.Range("J132").Value = 4128.7144
$ws.Range("K132").Value = 5465.9001
$ws.Range("L132").Value = 12386.1432
$ws.Range("M132").Value = -2935.9001
$ws.Range("N132").Value = -17446.1432
